$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.436250925064087
$ws.Range("B1").Value = 3.284970045089722
$ws.Range("C1").Value = 4.420218467712402
$ws.Range("D1").Value = 2.040529251098633
$ws.Range("E1").Value = 1.162777423858643
